# Auto-applies the 2026-02-13 06:20 meteocat refresh: updated extraction
# timestamps (col E) plus the handful of measurements that shifted between
# the 05:4x and 06:1x/06:2x scrapes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-13 06:18:21'
$ws.Range("E3").Value = '2026-02-13 06:18:24'
$ws.Range("E4").Value = '2026-02-13 06:18:26'
$ws.Range("H4").Value = '''57%'
$ws.Range("J4").Value = '1002.7 hPa'
$ws.Range("N4").Value = '8.2 °C 5:54 TU'
$ws.Range("O4").Value = '9.7 °C'
$ws.Range("E5").Value = '2026-02-13 06:18:28'
$ws.Range("H5").Value = '''66%'
$ws.Range("E6").Value = '2026-02-13 06:18:31'
$ws.Range("H6").Value = '''62%'
$ws.Range("J6").Value = '1002.9 hPa'
$ws.Range("E7").Value = '2026-02-13 06:18:33'
$ws.Range("J7").Value = '1003.4 hPa'
$ws.Range("N7").Value = '13.7 °C 5:35 TU'
$ws.Range("E8").Value = '2026-02-13 06:18:36'
$ws.Range("J8").Value = '1003.3 hPa'
$ws.Range("E9").Value = '2026-02-13 06:18:38'
$ws.Range("H9").Value = '''66%'
$ws.Range("O9").Value = '8.4 °C'
$ws.Range("E10").Value = '2026-02-13 06:18:41'
$ws.Range("H10").Value = '''78%'
$ws.Range("E11").Value = '2026-02-13 06:18:43'
$ws.Range("O11").Value = '2.1 °C'
$ws.Range("E12").Value = '2026-02-13 06:18:46'
$ws.Range("H12").Value = '''75%'
$ws.Range("O12").Value = '8.0 °C'
$ws.Range("E13").Value = '2026-02-13 06:18:48'
$ws.Range("J13").Value = '1006.9 hPa'
$ws.Range("E14").Value = '2026-02-13 06:18:51'
$ws.Range("H14").Value = '''62%'
$ws.Range("O14").Value = '11.4 °C'
$ws.Range("E15").Value = '2026-02-13 06:18:53'
$ws.Range("H15").Value = '''68%'
$ws.Range("O15").Value = '8.6 °C'
$ws.Range("E16").Value = '2026-02-13 06:18:55'
$ws.Range("H16").Value = '''61%'
$ws.Range("E17").Value = '2026-02-13 06:18:58'
$ws.Range("O17").Value = '1.4 °C'
$ws.Range("E18").Value = '2026-02-13 06:19:00'
$ws.Range("H18").Value = '''72%'
$ws.Range("J18").Value = '1003.0 hPa'
$ws.Range("E19").Value = '2026-02-13 06:19:03'
$ws.Range("E20").Value = '2026-02-13 06:19:05'
$ws.Range("H20").Value = '''78%'
$ws.Range("I20").Value = '0.1 mm'
$ws.Range("E21").Value = '2026-02-13 06:19:08'
$ws.Range("H21").Value = '''78%'
$ws.Range("J21").Value = '1005.5 hPa'
$ws.Range("O21").Value = '2.3 °C'
$ws.Range("E22").Value = '2026-02-13 06:19:10'
$ws.Range("H22").Value = '''83%'
$ws.Range("I22").Value = '0.3 mm'
$ws.Range("E23").Value = '2026-02-13 06:19:12'
$ws.Range("L23").Value = '47.2 km/h - 119º 5:45 TU'
$ws.Range("O23").Value = '-3.7 °C'
$ws.Range("E24").Value = '2026-02-13 06:19:15'
$ws.Range("I24").Value = '0.3 mm'
$ws.Range("J24").Value = '1004.4 hPa'
$ws.Range("O24").Value = '6.8 °C'
$ws.Range("E25").Value = '2026-02-13 06:19:17'
$ws.Range("H25").Value = '''59%'
$ws.Range("E26").Value = '2026-02-13 06:19:20'
$ws.Range("J26").Value = '1003.4 hPa'
$ws.Range("O26").Value = '2.6 °C'
$ws.Range("E27").Value = '2026-02-13 06:19:22'
$ws.Range("H27").Value = '''61%'
$ws.Range("I27").Value = '0.1 mm'
$ws.Range("E28").Value = '2026-02-13 06:19:25'
$ws.Range("H28").Value = '''67%'
$ws.Range("J28").Value = '1003.5 hPa'
$ws.Range("E29").Value = '2026-02-13 06:19:27'
$ws.Range("H29").Value = '''85%'
$ws.Range("L29").Value = '26.6 km/h - 164º 5:08 TU'
$ws.Range("O29").Value = '10.3 °C'
$ws.Range("E30").Value = '2026-02-13 06:19:29'
$ws.Range("J30").Value = '1003.2 hPa'
$ws.Range("O30").Value = '7.7 °C'
$ws.Range("E31").Value = '2026-02-13 06:19:32'
$ws.Range("H31").Value = '''55%'
$ws.Range("J31").Value = '1002.2 hPa'
$ws.Range("N31").Value = '10.5 °C 5:59 TU'
$ws.Range("O31").Value = '11.6 °C'
$ws.Range("E32").Value = '2026-02-13 06:19:34'
$ws.Range("I32").Value = '0.1 mm'
$ws.Range("O32").Value = '5.0 °C'
$ws.Range("E33").Value = '2026-02-13 06:19:37'
$ws.Range("H33").Value = '''76%'
$ws.Range("J33").Value = '1005.6 hPa'
$ws.Range("O33").Value = '0.6 °C'
$ws.Range("E34").Value = '2026-02-13 06:19:39'
$ws.Range("H34").Value = '''55%'
$ws.Range("O34").Value = '-0.4 °C'
$ws.Range("E35").Value = '2026-02-13 06:19:42'
$ws.Range("J35").Value = '1005.0 hPa'
$ws.Range("E36").Value = '2026-02-13 06:19:44'
$ws.Range("H36").Value = '''61%'
$ws.Range("J36").Value = '1002.9 hPa'
$ws.Range("O36").Value = '11.2 °C'
$ws.Range("E37").Value = '2026-02-13 06:19:47'
$ws.Range("H37").Value = '''64%'
$ws.Range("J37").Value = '1005.0 hPa'
$ws.Range("O37").Value = '3.8 °C'
$ws.Range("E38").Value = '2026-02-13 06:19:49'
$ws.Range("H38").Value = '''54%'
$ws.Range("K38").Value = '-0.1 MJ/m2'
$ws.Range("O38").Value = '10.4 °C'
$ws.Range("E39").Value = '2026-02-13 06:19:51'
$ws.Range("H39").Value = '''51%'
$ws.Range("N39").Value = '-4.2 °C 5:56 TU'
$ws.Range("O39").Value = '-2.6 °C'
$ws.Range("E40").Value = '2026-02-13 06:19:54'
$ws.Range("H40").Value = '''93%'
$ws.Range("J40").Value = '1006.7 hPa'
$ws.Range("E41").Value = '2026-02-13 06:19:56'
$ws.Range("H41").Value = '''51%'
$ws.Range("J41").Value = '1003.7 hPa'
$ws.Range("N41").Value = '10.1 °C 5:59 TU'
$ws.Range("O41").Value = '12.9 °C'
$ws.Range("E42").Value = '2026-02-13 06:19:58'
$ws.Range("O42").Value = '10.7 °C'
$ws.Range("E43").Value = '2026-02-13 06:20:01'
$ws.Range("H43").Value = '''66%'
$ws.Range("O43").Value = '7.2 °C'
$ws.Range("E44").Value = '2026-02-13 06:20:03'
$ws.Range("H44").Value = '''81%'
$ws.Range("L44").Value = '71.6 km/h - 199º 5:50 TU'
$ws.Range("E45").Value = '2026-02-13 06:20:06'
$ws.Range("H45").Value = '''68%'
$ws.Range("J45").Value = '1003.6 hPa'
$ws.Range("L45").Value = '19.1 km/h - 71º 5:39 TU'
$ws.Range("M45").Value = '5.5 °C 5:52 TU'
$ws.Range("O45").Value = '2.7 °C'
$ws.Range("E46").Value = '2026-02-13 06:20:08'
$ws.Range("H46").Value = '''82%'
$ws.Range("J46").Value = '1004.7 hPa'
$ws.Range("K46").Value = '-0.1 MJ/m2'
$ws.Range("O46").Value = '7.0 °C'
